{"js": "// Locate the single table in the document body.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\n// --- 1) Update the text of existing cells (indices based on the ---\n// --- original, pre-insert/pre-delete row order). ---------------\nconst values = table.values;\n\nvalues[0][0] = \"0M\";        // was 99.93\nvalues[1][0] = \"0M\";        // was 0.22\nvalues[2][0] = \"0M\";        // was 291\nvalues[3][0] = \"1208\";      // was 302\nvalues[5][0] = \"0.00254\";   // was 0.00242\nvalues[7][0] = \"0.00027\";   // was 0.00014\nvalues[8][0] = \"0.00030\";   // was 0.00007\nvalues[9][0] = \"0.00032\";   // was 0.00008\nvalues[10][0] = \"0.21569\";  // was 0.00010\n\n// The last three rows held tab-separated summary lines; they collapse\n// down to a single short value each.\nvalues[43][0] = \"99.93\";\nvalues[44][0] = \"0.22\";\nvalues[45][0] = \"291\";\n\ntable.values = values;\nawait context.sync();\n\n// --- 2) Insert a brand-new row (value \"0.00018\") right after the ---\n// --- row that holds \"0.00254\" (originally the \"0.00242\" row). ------\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst anchorRow = table.rows.items[5];\nanchorRow.insertRows(\"After\", 1, [[\"0.00018\"]]);\nawait context.sync();\n\n// --- 3) Delete the row that used to hold \"0.02753\". After the ------\n// --- insertion above, that row shifted from index 11 to index 12. --\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst rowToDelete = table.rows.items[12];\nrowToDelete.delete();\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# --- 1) Update the text of existing cells first (1-based row ---\n# --- indices, matching the table's original row order). --------\n$tbl.Cell(1, 1).Range.Text = \"0M\"        # was 99.93\n$tbl.Cell(2, 1).Range.Text = \"0M\"        # was 0.22\n$tbl.Cell(3, 1).Range.Text = \"0M\"        # was 291\n$tbl.Cell(4, 1).Range.Text = \"1208\"      # was 302\n$tbl.Cell(6, 1).Range.Text = \"0.00254\"   # was 0.00242\n$tbl.Cell(8, 1).Range.Text = \"0.00027\"   # was 0.00014\n$tbl.Cell(9, 1).Range.Text = \"0.00030\"   # was 0.00007\n$tbl.Cell(10, 1).Range.Text = \"0.00032\"  # was 0.00008\n$tbl.Cell(11, 1).Range.Text = \"0.21569\"  # was 0.00010\n\n# The last three rows held tab-separated summary lines; they collapse\n# down to a single short value each.\n$tbl.Cell(44, 1).Range.Text = \"99.93\"\n$tbl.Cell(45, 1).Range.Text = \"0.22\"\n$tbl.Cell(46, 1).Range.Text = \"291\"\n\n# --- 2) Insert a brand-new row (value \"0.00018\") right after the ---\n# --- row that holds \"0.00254\" (row 6, originally \"0.00242\"). -------\n$beforeRow = $tbl.Rows.Item(7)\n$newRow = $tbl.Rows.Add($beforeRow)\n$newRow.Cells.Item(1).Range.Text = \"0.00018\"\n\n# --- 3) Delete the row that used to hold \"0.02753\". After the ------\n# --- insertion above, that row shifted from row 12 to row 13. ------\n$tbl.Rows.Item(13).Delete()\n"}
